$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.529.50"

$ws.Range("D3").Value = "2.650.55"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "591.39"
$ws.Range("E5").Value = "  +1.84%  "

$ws.Range("D6").Value = "144.15"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").Value = "2.649.44"
$ws.Range("E9").Value = "  +2.89%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("E12").Value = "  +0.82%  "

$ws.Range("E13").Value = "  +0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").Value = "3.124.79"
$ws.Range("E15").Value = "  +2.93%  "

$ws.Range("D16").Value = "63.443.13"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").Value = "2.628.47"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").Value = "11.43"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").Value = "340.47"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("E22").Value = "  +1.50%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "1.67"
$ws.Range("E25").Value = "  +6.18%  "

$ws.Range("E26").Value = "  +5.05%  "

$ws.Range("E27").Value = "  +0.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "543.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.74%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +2.85%  "

$ws.Range("D31").Value = "7.78"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("E32").Value = "  +14.50%  "

$ws.Range("E33").Value = "  +2.72%  "

$ws.Range("D34").Value = "0.0₃0807"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "175.42"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").Value = "4.89"
$ws.Range("E36").Value = "  +9.59%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("D39").Value = "19.09"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("E40").Value = "  +7.79%  "

# Row 41 (was USDe) becomes Aave
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "170.77"
$ws.Range("E41").Value = "  +8.72%  "

# Row 42 (was Aave) becomes USDe
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "40.29"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").Value = "3.75"
$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("D45").Value = "22.45"
$ws.Range("E45").Value = "  +6.79%  "

$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").Value = "0.0558"
$ws.Range("E47").Value = "  +4.66%  "

$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("D50").Value = "18.86"
$ws.Range("E50").Value = "  +4.41%  "

$ws.Range("E51").Value = "  +0.75%  "
